# Update LR-pair TPM data: replace the 3-row (Sending x Ligand x Receptor)
# summary with the full 6-row breakdown per Sending/Target cluster pair,
# using the newly computed TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Myoc", "Fzd10", "ECs",   2, 0.6666666666666666, 0.113818,   0.341454,   0.0003230180320166274, 0.0003230180320166274, 1, 0.3333333333333333, 0.003092666666666667, 0.009278, 0.03934390080485799, 0.03934390080485798, 0.0003520011346666667, 0.003168010212,     0.00001270878940984263, 0.00001270878940984263),
    @("ECs",  "Myoc", "Fzd10", "MuSCs", 2, 0.6666666666666666, 0.113818,   0.341454,   0.0003230180320166274, 0.0003230180320166274, 2, 0.6666666666666666, 0.07551333333333334,  0.22654,  0.9606560991951421,  0.9606560991951419,  0.008594776573333332,  0.07735298916,      0.0003103092426067847, 0.0003103092426067847),
    @("FAPs", "Myoc", "Fzd10", "ECs",   3, 1,                  348.977468, 1046.932404, 0.9904058666599795,    0.9904058666599794,    1, 0.3333333333333333, 0.003092666666666667, 0.009278, 0.03934390080485799, 0.03934390080485798, 1.079270982701333,     9.713438844311998,  0.03896643017441964,   0.03896643017441963),
    @("FAPs", "Myoc", "Fzd10", "MuSCs", 3, 1,                  348.977468, 1046.932404, 0.9904058666599795,    0.9904058666599794,    2, 0.6666666666666666, 0.07551333333333334,  0.22654,  0.9606560991951421,  0.9606560991951419,  26.35245186690667,     237.17206680216,    0.95143943648556,      0.9514394364855597),
    @("MuSCs","Myoc", "Fzd10", "ECs",   3, 1,                  3.266752,   9.800255999999999, 0.009271115308003845, 0.009271115308003843, 1, 0.3333333333333333, 0.003092666666666667, 0.009278, 0.03934390080485799, 0.03934390080485798, 0.01010297501866667, 0.09092677516799999, 0.0003647618410285037, 0.0003647618410285035),
    @("MuSCs","Myoc", "Fzd10", "MuSCs", 3, 1,                  3.266752,   9.800255999999999, 0.009271115308003845, 0.009271115308003843, 2, 0.6666666666666666, 0.07551333333333334,  0.22654,  0.9606560991951421,  0.9606560991951419,  0.2466833326933333,  2.22014999424,       0.008906353466975341,  0.008906353466975339)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}
